# "Actualizar" run — 02-05-2021 03-44-35
# 1) re-stamp the previous check's timestamp block (rows 310:323)
# 2) append a brand-new check block (rows 324:337) with its own timestamp + hyperlinks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) touch the last existing block's Fecha column -----------------------
$refreshedStamp = 44232.13486850695
for ($r = 310; $r -le 323; $r++) {
    $ws.Range("D$r").Value = $refreshedStamp
}

# ---- 2) append the new block -------------------------------------------
$newStamp = 44232.15592659322

$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$linkSubAddresses = @("","","","","","","","","/","","","","","")

$startRow = 324
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i

    $ws.Range("A$row").Value = $names[$i]
    $ws.Range("C$row").Value = "Disponible"

    $ws.Range("D$row").Value = $newStamp
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B$row").Value = $displayUrls[$i]
    $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i], $linkSubAddresses[$i])
    $ws.Range("B$row").Style = "Hyperlink"
}

Write-Output "Updated D310:D323 and appended rows 324:337 with hyperlinks"
